$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A87").Value = "Optim_tuning"
$ws.Range("C87").Value = @'
001 Rf mapping start at 14:34
-8:2:8 2deg
'@
$ws.Range("D87").Value = "Alfa64chan-16032020-001"
$ws.Range("E87").Value = "200316_Alfa_rfMapper_basic"
$ws.Range("H87").Value = @'
N:\Stimuli\2019-06-RF-mapping\2020-03-16-Alfa
'@

$ws.Range("A88").Value = "Optim_tuning"
$ws.Range("C88").Value = @'
003 RF mapping starts 14:45
-1.5:0.75:1.5  0.75 deg
'@
$ws.Range("D88").Value = "Alfa64chan-16032020-003"
$ws.Range("E88").Value = "200316_Alfa_rfMapper_basic(2)"
$ws.Range("H88").Value = @'
N:\Stimuli\2019-06-RF-mapping\2020-03-16-Alfa
'@

$ws.Range("A89").Value = "Optim_tuning"
$ws.Range("C89").Value = @'
004 Generate Integrated 14:54
Testing if we don't use baseline subtracted score, can we evolve faster.
12 [-0.5 -0.8] 1  3 deg
30 gens
Finish in 16 mins  plateau around 16 gens! Grows super fast!
Very successful!
'@
$ws.Range("D89").Value = "Alfa64chan-16032020-004"
$ws.Range("E89").Value = "200316_Alfa_generate_integrated(1)"
$ws.Range("H89").Value = @'
N:\Stimuli\2019-12-Evolutions\2020-03-16-Alfa-01\2020-03-16-14-54-33
'@

$ws.Range("A90").Value = "Optim_tuning"
$ws.Range("C90").Value = @'
005 Generate integrated 15:14
12 [-0.5 -0.8] 1  3 deg   ZOHA_Sphere_lr_euclid
Mu exponential decay from 50 deg to 7.33deg in 100 generations, learning rate 1.5
Try out my new optimizer!
Optimizer found gradient around 16-18 gens , grow featues to excite the plateau part of the image.
Exploration get down to 32 deg at gen 24.
Seems to surpass the performance of CMAES finally!
'@
$ws.Range("D90").Value = "Alfa64chan-16032020-005"
$ws.Range("E90").Value = "200316_Alfa_generate_integrated(2)"
$ws.Range("H90").Value = @'
N:\Stimuli\2019-12-Evolutions\2020-03-16-Alfa-02\2020-03-16-15-14-25
'@

$ws.Range("A91").Value = "Optim_tuning"
$ws.Range("C91").Value = @'
006 Generate integrated 15:40
12 [-0.5 -0.8] 1  3 deg   ZOHA_Sphere_lr_euclid
Mu exponential decay from 50 deg to 20 deg in 100 generations, learning rate 1.5
See if larger exploration will result in different evolutions trajectory?
Want to test if I decayed the learning rate too fast or too slow, because both could result in a slow optimization.
See which part of the experiment does the firing rate grows the fastest! May be that is the right rate to go.
Starts growing around gen 16-18, exploration = 41.5. At gen 23 the step between basis are 0.442.
Close to 30 gens, still around 38.9 degree exploration.
Seems the plateau around 150ms doesn't evolve successfully
Starts growing the plateau around 32 gens.
Obviously, the evolution is much less smooth than the last one! Many stagnant and coming in and out.
Around 35 gens still 36.5 degrees.
37 gens 35.8 degree,
This is much slower than the 2nd, ZOHA Sphere_lr_euclide, and slower than the first CMA-ES,
PSTH seems different from the last one, the gap between the 2 peaks are longer!
At around 50 gens reached 450 just
This schedule seems not very helpful……. Try others
'@
$ws.Range("D91").Value = "Alfa64chan-16032020-006"
$ws.Range("E91").Value = "200316_Alfa_generate_integrated(3)"
$ws.Range("H91").Value = @'
N:\Stimuli\2019-12-Evolutions\2020-03-16-Alfa-03\2020-03-16-15-40-58
'@

$ws.Range("A92").Value = "Optim_tuning"
$ws.Range("C92").Value = @'
007 Generate Integrated 16:12
12 [-0.5 -0.8] 1  3 deg   ZOHA_Sphere_lr_euclid
Mu exponential decay from 40 deg to 10 deg in 100 generations, learning rate 1.5
See if smaller exploration will result in better evolutions trajectory?
At gen 9, the exploration around 35.3 deg
Speed of this optimizer is similar to CMAES! Grows super fast, Exploration range around 32.9 at 15gens.
Perhaps exploration range around 30-35deg is the sweet spot?
Reach 400 around 30 gens, similar PSTH.
Faster than the last one. But curiously, the result doesn't look as
'@
$ws.Range("D92").Value = "Alfa64chan-16032020-007"
$ws.Range("E92").Value = "200316_Alfa_generate_integrated(4)"
$ws.Range("H92").Value = @'
N:\Stimuli\2019-12-Evolutions\2020-03-16-Alfa-04\2020-03-16-16-12-37
'@

$ws.Range("A93").Value = "Optim_tuning"
$ws.Range("C93").Value = @'
008 Generate Integrated 16:45
64 [0 -2.1] 1  4deg  ZOHA_Sphere_lr_euclid
Same learning rate scheduling with last experiment,
Mu exponential decay from 40 deg to 10 deg in 100 generations, learning rate 1.5
The firing rate PSTH is very sustained! Like square wave.
Starts growing steadily from 5 gens.
Grow not very fast. Starts to plateau around 20 gens.
But the learning rate schedule tuned for
'@
$ws.Range("D93").Value = "Alfa64chan-16032020-008"
$ws.Range("E93").Value = "200316_Alfa_generate_integrated(5)"
$ws.Range("H93").Value = @'
N:\Stimuli\2019-12-Evolutions\2020-03-16-Alfa-05\2020-03-16-16-45-43
'@

$ws.Range("A94").Value = "Optim_tuning"
$ws.Range("C94").Value = @'
009 Generate Integrated 17:03
64 [0 -2.1]  1 4deg  ZOHA_Sphere_lr_euclid
Mu exponential decay from 40 deg to 5 deg in 100 generations, learning rate 1.5
Faster decay to small learning rate, see how it goes!
The result doesn't change significantly! Seems similar speed and result
He is a very good boy today! Seems he can give me 20 gens for this
'@
$ws.Range("D94").Value = "Alfa64chan-16032020-009"
$ws.Range("E94").Value = "200316_Alfa_generate_integrated(6)"
$ws.Range("H94").Value = @'
N:\Stimuli\2019-12-Evolutions\2020-03-16-Alfa-06\2020-03-16-17-03-32
'@

$ws.Range("87:94").AutoFit()

